$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 4.570121
$ws.Range("D3").Value = 4.429555
$ws.Range("D4").Value = 4.3924713
$ws.Range("D5").Value = 4.475478
$ws.Range("D6").Value = 4.8108716
$ws.Range("D7").Value = 4.62121
$ws.Range("D8").Value = 4.427797
$ws.Range("D9").Value = 4.4477773
$ws.Range("D10").Value = 4.4033685
$ws.Range("D11").Value = 20.044765
$ws.Range("D12").Value = 20.09938
$ws.Range("D13").Value = 20.093447
$ws.Range("D14").Value = 19.954018
$ws.Range("D15").Value = 20.11724
$ws.Range("D16").Value = 20.191637
$ws.Range("D17").Value = 19.834846
$ws.Range("D18").Value = 19.86631
$ws.Range("D19").Value = 19.858034
$ws.Range("D20").Value = 5.340901
$ws.Range("D21").Value = 5.433157
$ws.Range("D22").Value = 5.1576967
$ws.Range("D23").Value = 4.4378176
$ws.Range("D24").Value = 4.4187546
$ws.Range("D25").Value = 4.51871
$ws.Range("D26").Value = 4.3101935
$ws.Range("D27").Value = 4.423321
$ws.Range("D28").Value = 4.4769
$ws.Range("D29").Value = 21.09741
$ws.Range("D30").Value = 20.952942
$ws.Range("D31").Value = 20.861547
$ws.Range("D32").Value = 20.21614
$ws.Range("D33").Value = 20.162424
$ws.Range("D34").Value = 20.3171
$ws.Range("D35").Value = 19.750463
$ws.Range("D36").Value = 19.930487
$ws.Range("D37").Value = 19.97378
$ws.Range("D38").Value = 4.8195386
$ws.Range("D39").Value = 4.5988398
$ws.Range("D40").Value = 4.704048
$ws.Range("D41").Value = 4.4755135
$ws.Range("D42").Value = 4.471478
$ws.Range("D43").Value = 4.465045
$ws.Range("D44").Value = 4.449566
$ws.Range("D45").Value = 4.3526654
$ws.Range("D46").Value = 4.479339
$ws.Range("D47").Value = 20.194372
$ws.Range("D48").Value = 20.077696
$ws.Range("D49").Value = 20.025795
$ws.Range("D50").Value = 20.066381
$ws.Range("D51").Value = 19.973536
$ws.Range("D52").Value = 20.033766
$ws.Range("D53").Value = 19.903559
$ws.Range("D54").Value = 19.973228
$ws.Range("D55").Value = 19.778542
$ws.Range("D56").Value = 4.989229
$ws.Range("D57").Value = 4.570296
$ws.Range("D58").Value = 4.520977
$ws.Range("D59").Value = 4.466665
$ws.Range("D60").Value = 4.3805895
$ws.Range("D61").Value = 4.465386
$ws.Range("D62").Value = 4.540096
$ws.Range("D63").Value = 4.4209065
$ws.Range("D64").Value = 4.407012
$ws.Range("D65").Value = 20.196901
$ws.Range("D66").Value = 20.066807
$ws.Range("D67").Value = 20.043415
$ws.Range("D68").Value = 19.957195
$ws.Range("D69").Value = 19.960686
$ws.Range("D70").Value = 19.969162
$ws.Range("D71").Value = 20.160696
$ws.Range("D72").Value = 20.111479
$ws.Range("D73").Value = 20.14878
